$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (index 1 in the data table) now corresponds to the 380 kV case:
# B2..E2 go from 0 -> the computed load/p_mw results.
$ws.Range("B2").Value = 256.5101652550508
$ws.Range("C2").Value = 95.76379502855229
$ws.Range("D2").Value = 95.76379502855229
$ws.Range("E2").Value = 95.76379502855229
